# Re-style the three tables that currently use the default "Table_0"
# style ({F5EECA30-BCA4-4A58-A9FE-A3DDC721C4EE}) so that they instead use
# the "Medium Style 2 - Accent 1" built-in style ({3A972697-7673-4304-9457-066B01AFF1A1}).
#
# PowerPoint's Table object does not let you assign .Style directly
# ("Table styles cannot be assigned through a property ... call
# Table.ApplyStyle(\"{GUID}\") instead"), so ApplyStyle is used.

$OldStyleId = "{F5EECA30-BCA4-4A58-A9FE-A3DDC721C4EE}"
$NewStyleId = "{3A972697-7673-4304-9457-066B01AFF1A1}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $OldStyleId) {
                $table.ApplyStyle($NewStyleId)
            }
        }
    }
}
